$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - name unchanged, values updated
$ws.Range("B3").Value = 0.9857675710868646
$ws.Range("C3").Value = 0.985351507780873
$ws.Range("D3").Value = 0.9846258810617178

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.984882949614016
$ws.Range("C4").Value = 0.9852420679042488
$ws.Range("D4").Value = 0.9849952715470319

# Row 5: AdaBoostRegressor -> MLPRegressor, values updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9893350389276304
$ws.Range("C5").Value = 0.9890726811129923
$ws.Range("D5").Value = 0.9894529453940709
